$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: change absolute $A$2 reference to relative A2 in the formulas
# for columns C, D, E (matching the pattern already used in column B).
$ws.Range("C2").Formula = '=$C$1 & TEXT(A2, "##")'
$ws.Range("D2").Formula = '=$D$1 & TEXT(A2, "##")'
$ws.Range("E2").Formula = '=$E$1 & TEXT(A2, "##")'

# Rows 3:17: columns D and E were using a shared formula that still
# referenced the absolute $A$2 cell, so every row evaluated to the same
# "Marks1"/"Gender1" text. Re-apply the formula across the range so it
# mirrors columns B/C (relative reference to the row's own A cell) and
# restores the shared-formula behaviour with correctly recalculated values.
$ws.Range("D3:D17").Formula = '=$D$1 & TEXT(A3, "##")'
$ws.Range("E3:E17").Formula = '=$E$1 & TEXT(A3, "##")'
